# Weekly Fruit/Vegetable price update: insert a new daily record for
# "Macroferia Regional de Talca" (Mango) as row 49, pushing the existing
# rows 49-130 down by one (to 50-131).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 49 - shifts rows 49..130 down to 50..131
$ws.Rows(49).Insert()

# Populate the newly inserted row 49 with the new record's data
$ws.Range("A49").Value = 5
$ws.Range("B49").Value = "Macroferia Regional de Talca"
$ws.Range("C49").Value = "Maule"
$ws.Range("D49").Value = 44720
$ws.Range("E49").Value = 7
$ws.Range("F49").Value = "Fruta"
$ws.Range("G49").Value = 100108
$ws.Range("H49").Value = "Tropicales y subtropicales"
$ws.Range("I49").Value = 100108002
$ws.Range("J49").Value = "Mango"
$ws.Range("K49").Value = "Sin especificar"
$ws.Range("L49").Value = "Primera"
$ws.Range("M49").Value = 241
$ws.Range("N49").Value = 10000
$ws.Range("O49").Value = 10000
$ws.Range("P49").Value = 10000
$ws.Range("Q49").Value = "$/bandeja 4 kilos"
$ws.Range("R49").Value = "Brasil"
$ws.Range("S49").Value = 2500
$ws.Range("T49").Value = 4
